# Auto-generated edit script: updates crypto price/volume table cells
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.548.89"
$ws.Range("E2").Value = "  -1.02%  "
$ws.Range("D3").Value = "1.924.29"
$ws.Range("E3").Value = "  +1.78%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9991"
$ws.Range("D4").Style = "Normal"
$ws.Range("E5").Value = "  +1.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9993"
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = "  -1.89%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2883"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.36%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06770"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "106.56"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "18.35"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07755"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.41%  "
$ws.Range("D13").Value = "1.912.42"
$ws.Range("E13").Value = "  +1.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.292"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6616"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "293.54"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.90%  "
$ws.Range("D17").Value = "30.555.41"
$ws.Range("E17").Value = "  -0.90%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007607"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.56%  "
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.97%  "
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9991"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").Value = "2.159.83"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.297"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.90%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9988"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.219"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.372"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.39"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.115"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1075"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.365"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.184"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.010"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05043"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7439"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.156"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02096"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.725"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.689"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.077"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "110.39"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8745"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.923"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.49%  "
$ws.Range("E43").Value = "  +0.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9992"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "67.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "50.09"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +18.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.226"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.323"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.24%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1221"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.10%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "35.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.2467"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +8.35%  "
